$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Phase 1 Pre CPP")

# Delete rows 3 through 9 (the extra data rows being removed)
$ws.Range("A3:F9").EntireRow.Delete() | Out-Null

# Update the remaining data row (row 2) with corrected values
$ws.Range("A2").Value = 30414
$ws.Range("B2").Value = 0.035
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 2061400
$ws.Range("E2").Value = 260.8258482527681
$ws.Range("F2").Value = 0.008641033266058135
